$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-08 (row 21)
$ws.Range("B21").Value = 6235
$ws.Range("C21").Value = 987
$ws.Range("D21").Value = 5623375
$ws.Range("E21").Value = 901.9045709703288
$ws.Range("F21").Value = 8.227738239888915
$ws.Range("G21").Value = 4.113924050632911
$ws.Range("H21").Value = 28.35907191372142
